# CO12AL-W2-VIDEO06-SLIDE01 : recolor the Python code sample to match the
# IDLE editor color scheme (orange keywords, green strings).
#
# Slide 1 holds a single content placeholder with the 6-line code sample:
#   1: note = 8
#   2: if note > 10:
#   3:     print 'reçu'
#   4:     print 'bravo !'
#   5: else:
#   6:     print 'recalé'
#
# We keep the text itself unchanged and only colorize sub-strings of each
# line: Python keywords (if / print / else) become orange (FFC000) and the
# quoted string literals become green (00B050). Using TextRange.Characters
# on each paragraph automatically splits the existing run and preserves all
# other run-level formatting (font, err spell-flag, etc.), matching how
# PowerPoint itself performs an in-place partial recolor.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Font.Color.RGB uses the classic OLE "BGR" integer (R + G*256 + B*65536),
# i.e. the reverse byte order of the hex srgbClr value used in the XML.
$ORANGE = 49407    # 0xFFC000 -> srgbClr val="FFC000"
$GREEN  = 5287936  # 0x00B050 -> srgbClr val="00B050"

# Paragraph 2: "if note > 10:"  ->  "if" in orange
$para = $tr.Paragraphs(2)
$para.Characters(1, 2).Font.Color.RGB = $ORANGE

# Paragraph 3: "    print 'reçu'"  -> "print" orange, "'reçu'" green
$para = $tr.Paragraphs(3)
$para.Characters(5, 5).Font.Color.RGB = $ORANGE
$para.Characters(11, 6).Font.Color.RGB = $GREEN

# Paragraph 4: "    print 'bravo !'" -> "print" orange, "'bravo !'" green
$para = $tr.Paragraphs(4)
$para.Characters(5, 5).Font.Color.RGB = $ORANGE
$para.Characters(11, 9).Font.Color.RGB = $GREEN

# Paragraph 5: "else:" -> "else" orange
$para = $tr.Paragraphs(5)
$para.Characters(1, 4).Font.Color.RGB = $ORANGE

# Paragraph 6: "    print 'recalé'" -> "print" orange, "'recalé'" green
$para = $tr.Paragraphs(6)
$para.Characters(5, 5).Font.Color.RGB = $ORANGE
$para.Characters(11, 8).Font.Color.RGB = $GREEN
